$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 8).Value = 1172.2354  # H17: 1313.1034 -> 1172.2354
$ws.Cells.Item(17, 10).Value = 1172.2354  # J17: 1313.1034 -> 1172.2354
$ws.Cells.Item(17, 12).Value = 3516.7062  # L17: 3939.3102 -> 3516.7062
$ws.Cells.Item(17, 14).Value = -3852.7062  # N17: -4275.3102 -> -3852.7062
$ws.Cells.Item(19, 8).Value = 1715  # H19: 1939.3846 -> 1715
$ws.Cells.Item(19, 9).Value = 420  # I19: 487.25 -> 420
$ws.Cells.Item(19, 10).Value = 2362.5  # J19: 2584.7778 -> 2362.5
$ws.Cells.Item(19, 11).Value = 420  # K19: 487.25 -> 420
$ws.Cells.Item(19, 12).Value = 2362.5  # L19: 2584.7778 -> 2362.5
$ws.Cells.Item(19, 13).Value = -245  # M19: -312.25 -> -245
$ws.Cells.Item(19, 14).Value = -2712.5  # N19: -2934.7778 -> -2712.5
$ws.Cells.Item(33, 8).Value = 62500108  # H33: 66666780 -> 62500108
$ws.Cells.Item(33, 9).Value = 118.46154  # I33: 122.583336 -> 118.46154
$ws.Cells.Item(33, 11).Value = 118.46154  # K33: 122.583336 -> 118.46154
$ws.Cells.Item(33, 13).Value = 110.53846  # M33: 106.416664 -> 110.53846
$ws.Cells.Item(43, 8).Value = 3500.4285  # H43: 3111.3333 -> 3500.4285
$ws.Cells.Item(43, 9).Value = 2935  # I43: 2501 -> 2935
$ws.Cells.Item(43, 10).Value = 3924.5  # J43: 3599.6 -> 3924.5
$ws.Cells.Item(43, 11).Value = 2935  # K43: 2501 -> 2935
$ws.Cells.Item(43, 12).Value = 3924.5  # L43: 3599.6 -> 3924.5
$ws.Cells.Item(43, 13).Value = -2866  # M43: -2432 -> -2866
$ws.Cells.Item(43, 14).Value = -4062.5  # N43: -3737.6 -> -4062.5
$ws.Cells.Item(54, 8).Value = 1209.8  # H54: 1210.2 -> 1209.8
$ws.Cells.Item(54, 9).Value = 1209.8  # I54: 1210.2 -> 1209.8
$ws.Cells.Item(54, 11).Value = 1209.8  # K54: 1210.2 -> 1209.8
$ws.Cells.Item(54, 13).Value = -723.8  # M54: -724.2 -> -723.8
$ws.Cells.Item(57, 8).Value = 38499.5  # H57: 46587.5 -> 38499.5
$ws.Cells.Item(57, 10).Value = 38499.5  # J57: 46587.5 -> 38499.5
$ws.Cells.Item(57, 12).Value = 115498.5  # L57: 139762.5 -> 115498.5
$ws.Cells.Item(57, 14).Value = -116496.5  # N57: -140760.5 -> -116496.5
$ws.Cells.Item(64, 8).Value = 3974.4  # H64: 3998.6667 -> 3974.4
$ws.Cells.Item(64, 10).Value = 3893.5  # J64: 3938.8 -> 3893.5
$ws.Cells.Item(64, 12).Value = 3893.5  # L64: 3938.8 -> 3893.5
$ws.Cells.Item(64, 14).Value = -4389.5  # N64: -4434.8 -> -4389.5
$ws.Cells.Item(67, 8).Value = 3974.4  # H67: 3998.6667 -> 3974.4
$ws.Cells.Item(67, 10).Value = 3893.5  # J67: 3938.8 -> 3893.5
$ws.Cells.Item(67, 12).Value = 3893.5  # L67: 3938.8 -> 3893.5
$ws.Cells.Item(67, 14).Value = -5609.5  # N67: -5654.8 -> -5609.5
$ws.Cells.Item(80, 8).Value = 818.25  # H80: 750.7 -> 818.25
$ws.Cells.Item(80, 9).Value = 526.7  # I80: 526.8 -> 526.7
$ws.Cells.Item(80, 10).Value = 1109.8  # J80: 974.6 -> 1109.8
$ws.Cells.Item(80, 11).Value = 1580.1  # K80: 1580.4 -> 1580.1
$ws.Cells.Item(80, 12).Value = 3329.4  # L80: 2923.8 -> 3329.4
$ws.Cells.Item(80, 13).Value = -582.1000000000001  # M80: -582.3999999999999 -> -582.1000000000001
$ws.Cells.Item(80, 14).Value = -5325.4  # N80: -4919.8 -> -5325.4
$ws.Cells.Item(83, 8).Value = 818.25  # H83: 750.7 -> 818.25
$ws.Cells.Item(83, 9).Value = 526.7  # I83: 526.8 -> 526.7
$ws.Cells.Item(83, 10).Value = 1109.8  # J83: 974.6 -> 1109.8
$ws.Cells.Item(83, 11).Value = 4740.3  # K83: 4741.2 -> 4740.3
$ws.Cells.Item(83, 12).Value = 9988.199999999999  # L83: 8771.4 -> 9988.199999999999
$ws.Cells.Item(83, 13).Value = 251.6999999999998  # M83: 250.8000000000002 -> 251.6999999999998
$ws.Cells.Item(83, 14).Value = -19972.2  # N83: -18755.4 -> -19972.2
$ws.Cells.Item(88, 8).Value = 717701.9399999999  # H88: 457031 -> 717701.9399999999
$ws.Cells.Item(88, 9).Value = 3250  # I88: 1792 -> 3250
$ws.Cells.Item(88, 10).Value = 836777.25  # J88: 558195.25 -> 836777.25
$ws.Cells.Item(88, 11).Value = 3250  # K88: 1792 -> 3250
$ws.Cells.Item(88, 12).Value = 836777.25  # L88: 558195.25 -> 836777.25
$ws.Cells.Item(88, 13).Value = -2844  # M88: -1386 -> -2844
$ws.Cells.Item(88, 14).Value = -837589.25  # N88: -559007.25 -> -837589.25
$ws.Cells.Item(91, 8).Value = 717701.9399999999  # H91: 457031 -> 717701.9399999999
$ws.Cells.Item(91, 9).Value = 3250  # I91: 1792 -> 3250
$ws.Cells.Item(91, 10).Value = 836777.25  # J91: 558195.25 -> 836777.25
$ws.Cells.Item(91, 11).Value = 3250  # K91: 1792 -> 3250
$ws.Cells.Item(91, 12).Value = 836777.25  # L91: 558195.25 -> 836777.25
$ws.Cells.Item(91, 13).Value = -1846  # M91: -388 -> -1846
$ws.Cells.Item(91, 14).Value = -839585.25  # N91: -561003.25 -> -839585.25
$ws.Cells.Item(116, 8).Value = 4443.1113  # H116: 4068 -> 4443.1113
$ws.Cells.Item(116, 9).Value = 4355.5713  # I116: 3963.1538 -> 4355.5713
$ws.Cells.Item(116, 11).Value = 4355.5713  # K116: 3963.1538 -> 4355.5713
$ws.Cells.Item(116, 13).Value = -913.5712999999996  # M116: -521.1538 -> -913.5712999999996
$ws.Cells.Item(132, 8).Value = 781.1905  # H132: 723.25 -> 781.1905
$ws.Cells.Item(132, 9).Value = 781.1905  # I132: 748.8889 -> 781.1905
$ws.Cells.Item(132, 10).Value = 0  # J132: 558.4286 -> 0
$ws.Cells.Item(132, 11).Value = 2343.5715  # K132: 2246.6667 -> 2343.5715
$ws.Cells.Item(132, 12).Value = 0  # L132: 1675.2858 -> 0
$ws.Cells.Item(132, 13).Value = 186.4285  # M132: 283.3332999999998 -> 186.4285
$ws.Cells.Item(132, 14).ClearContents()  # N132 was -6735.2858
$ws.Cells.Item(137, 8).Value = 1999.7931  # H137: 1865.1177 -> 1999.7931
$ws.Cells.Item(137, 9).Value = 2564.3333  # I137: 2718.4285 -> 2564.3333
$ws.Cells.Item(137, 10).Value = 1394.9286  # J137: 1267.8 -> 1394.9286
$ws.Cells.Item(137, 11).Value = 7692.999899999999  # K137: 8155.2855 -> 7692.999899999999
$ws.Cells.Item(137, 12).Value = 4184.7858  # L137: 3803.4 -> 4184.7858
$ws.Cells.Item(137, 13).Value = -5142.999899999999  # M137: -5605.2855 -> -5142.999899999999
$ws.Cells.Item(137, 14).Value = -9284.7858  # N137: -8903.4 -> -9284.7858
$ws.Cells.Item(138, 8).Value = 1929.1746  # H138: 2000.5862 -> 1929.1746
$ws.Cells.Item(138, 9).Value = 1076.9333  # I138: 1128.3846 -> 1076.9333
$ws.Cells.Item(138, 10).Value = 2195.5  # J138: 2252.5557 -> 2195.5
$ws.Cells.Item(138, 11).Value = 3230.7999  # K138: 3385.1538 -> 3230.7999
$ws.Cells.Item(138, 12).Value = 6586.5  # L138: 6757.6671 -> 6586.5
$ws.Cells.Item(138, 13).Value = 1909.2001  # M138: 1754.8462 -> 1909.2001
$ws.Cells.Item(138, 14).Value = -16866.5  # N138: -17037.6671 -> -16866.5
$ws.Cells.Item(141, 8).Value = 4012.8696  # H141: 4158.591 -> 4012.8696
$ws.Cells.Item(141, 9).Value = 3447.4285  # I141: 3579.45 -> 3447.4285
$ws.Cells.Item(141, 11).Value = 10342.2855  # K141: 10738.35 -> 10342.2855
$ws.Cells.Item(141, 13).Value = -5162.2855  # M141: -5558.349999999999 -> -5162.2855

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(45, 8).Value = 1820.5238  # H45: 1826.238 -> 1820.5238
$ws.Cells.Item(45, 9).Value = 1752.0667  # I45: 1818.2142 -> 1752.0667
$ws.Cells.Item(45, 10).Value = 1991.6666  # J45: 1842.2858 -> 1991.6666
$ws.Cells.Item(45, 11).Value = 1752.0667  # K45: 1818.2142 -> 1752.0667
$ws.Cells.Item(45, 12).Value = 1991.6666  # L45: 1842.2858 -> 1991.6666
$ws.Cells.Item(45, 13).Value = -1375.0667  # M45: -1441.2142 -> -1375.0667
$ws.Cells.Item(45, 14).Value = -2745.6666  # N45: -2596.2858 -> -2745.6666
$ws.Cells.Item(61, 8).Value = 1555.1765  # H61: 1730.9231 -> 1555.1765
$ws.Cells.Item(61, 9).Value = 1528.3077  # I61: 1572.909 -> 1528.3077
$ws.Cells.Item(61, 10).Value = 1642.5  # J61: 2600 -> 1642.5
$ws.Cells.Item(61, 11).Value = 1528.3077  # K61: 1572.909 -> 1528.3077
$ws.Cells.Item(61, 12).Value = 1642.5  # L61: 2600 -> 1642.5
$ws.Cells.Item(61, 13).Value = -1316.3077  # M61: -1360.909 -> -1316.3077
$ws.Cells.Item(61, 14).Value = -2066.5  # N61: -3024 -> -2066.5
$ws.Cells.Item(74, 8).Value = 2360.7896  # H74: 2545.9429 -> 2360.7896
$ws.Cells.Item(74, 9).Value = 1831.24  # I74: 2053.5908 -> 1831.24
$ws.Cells.Item(74, 11).Value = 1831.24  # K74: 2053.5908 -> 1831.24
$ws.Cells.Item(74, 13).Value = -957.24  # M74: -1179.5908 -> -957.24
$ws.Cells.Item(77, 8).Value = 2360.7896  # H77: 2545.9429 -> 2360.7896
$ws.Cells.Item(77, 9).Value = 1831.24  # I77: 2053.5908 -> 1831.24
$ws.Cells.Item(77, 11).Value = 9156.200000000001  # K77: 10267.954 -> 9156.200000000001
$ws.Cells.Item(77, 13).Value = -4788.200000000001  # M77: -5899.954 -> -4788.200000000001
$ws.Cells.Item(110, 8).Value = 1193.625  # H110: 1258 -> 1193.625
$ws.Cells.Item(110, 9).Value = 1193.625  # I110: 1258 -> 1193.625
$ws.Cells.Item(110, 11).Value = 1193.625  # K110: 1258 -> 1193.625
$ws.Cells.Item(110, 13).Value = 851.375  # M110: 787 -> 851.375
$ws.Cells.Item(122, 8).Value = 4070.5334  # H122: 4311.186 -> 4070.5334
$ws.Cells.Item(122, 9).Value = 3942.4722  # I122: 4094 -> 3942.4722
$ws.Cells.Item(122, 10).Value = 4582.778  # J122: 5428.143 -> 4582.778
$ws.Cells.Item(122, 11).Value = 11827.4166  # K122: 12282 -> 11827.4166
$ws.Cells.Item(122, 12).Value = 13748.334  # L122: 16284.429 -> 13748.334
$ws.Cells.Item(122, 13).Value = -9377.4166  # M122: -9832 -> -9377.4166
$ws.Cells.Item(122, 14).Value = -18648.334  # N122: -21184.429 -> -18648.334
$ws.Cells.Item(132, 8).Value = 11859.182  # H132: 7186.684 -> 11859.182
$ws.Cells.Item(132, 9).Value = 7544.85  # I132: 4631.1143 -> 7544.85
$ws.Cells.Item(132, 10).Value = 55002.5  # J132: 37001.668 -> 55002.5
$ws.Cells.Item(132, 11).Value = 22634.55  # K132: 13893.3429 -> 22634.55
$ws.Cells.Item(132, 12).Value = 165007.5  # L132: 111005.004 -> 165007.5
$ws.Cells.Item(132, 13).Value = -20104.55  # M132: -11363.3429 -> -20104.55
$ws.Cells.Item(132, 14).Value = -170067.5  # N132: -116065.004 -> -170067.5
$ws.Cells.Item(136, 8).Value = 1555.1765  # H136: 1730.9231 -> 1555.1765
$ws.Cells.Item(136, 9).Value = 1528.3077  # I136: 1572.909 -> 1528.3077
$ws.Cells.Item(136, 10).Value = 1642.5  # J136: 2600 -> 1642.5
$ws.Cells.Item(136, 11).Value = 4584.9231  # K136: 4718.727000000001 -> 4584.9231
$ws.Cells.Item(136, 12).Value = 4927.5  # L136: 7800 -> 4927.5
$ws.Cells.Item(136, 13).Value = -2034.9231  # M136: -2168.727000000001 -> -2034.9231
$ws.Cells.Item(136, 14).Value = -10027.5  # N136: -12900 -> -10027.5

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 1071.275  # H20: 1049.4054 -> 1071.275
$ws.Cells.Item(20, 9).Value = 1142.4073  # I20: 1117.5834 -> 1142.4073
$ws.Cells.Item(20, 11).Value = 1142.4073  # K20: 1117.5834 -> 1142.4073
$ws.Cells.Item(20, 13).Value = -895.4073000000001  # M20: -870.5834 -> -895.4073000000001
$ws.Cells.Item(86, 8).Value = 2476.0908  # H86: 2480.7727 -> 2476.0908
$ws.Cells.Item(86, 9).Value = 2311.2856  # I86: 2318.6428 -> 2311.2856
$ws.Cells.Item(86, 11).Value = 2311.2856  # K86: 2318.6428 -> 2311.2856
$ws.Cells.Item(86, 13).Value = -1188.2856  # M86: -1195.6428 -> -1188.2856
$ws.Cells.Item(89, 8).Value = 2476.0908  # H89: 2480.7727 -> 2476.0908
$ws.Cells.Item(89, 9).Value = 2311.2856  # I89: 2318.6428 -> 2311.2856
$ws.Cells.Item(89, 11).Value = 11556.428  # K89: 11593.214 -> 11556.428
$ws.Cells.Item(89, 13).Value = -5940.428  # M89: -5977.214 -> -5940.428
$ws.Cells.Item(94, 8).Value = 4253.647  # H94: 5477.0713 -> 4253.647
$ws.Cells.Item(94, 9).Value = 3111.8  # I94: 3423.2222 -> 3111.8
$ws.Cells.Item(94, 10).Value = 5884.857  # J94: 9174 -> 5884.857
$ws.Cells.Item(94, 11).Value = 3111.8  # K94: 3423.2222 -> 3111.8
$ws.Cells.Item(94, 12).Value = 5884.857  # L94: 9174 -> 5884.857
$ws.Cells.Item(94, 13).Value = -2660.8  # M94: -2972.2222 -> -2660.8
$ws.Cells.Item(94, 14).Value = -6786.857  # N94: -10076 -> -6786.857
$ws.Cells.Item(105, 8).Value = 2881.8  # H105: 3071.8333 -> 2881.8
$ws.Cells.Item(105, 9).Value = 2565.8  # I105: 2614 -> 2565.8
$ws.Cells.Item(105, 10).Value = 3197.8  # J105: 3529.6667 -> 3197.8
$ws.Cells.Item(105, 11).Value = 2565.8  # K105: 2614 -> 2565.8
$ws.Cells.Item(105, 12).Value = 3197.8  # L105: 3529.6667 -> 3197.8
$ws.Cells.Item(105, 13).Value = -818.8000000000002  # M105: -867 -> -818.8000000000002
$ws.Cells.Item(105, 14).Value = -6691.8  # N105: -7023.6667 -> -6691.8
$ws.Cells.Item(134, 8).Value = 2649.2666  # H134: 2454.1738 -> 2649.2666
$ws.Cells.Item(134, 9).Value = 2435.4167  # I134: 2187.1428 -> 2435.4167
$ws.Cells.Item(134, 10).Value = 3504.6667  # J134: 3303.818 -> 3504.6667
$ws.Cells.Item(134, 11).Value = 7306.250100000001  # K134: 6561.428400000001 -> 7306.250100000001
$ws.Cells.Item(134, 12).Value = 10514.0001  # L134: 9911.454000000002 -> 10514.0001
$ws.Cells.Item(134, 13).Value = -4771.250100000001  # M134: -4026.428400000001 -> -4771.250100000001
$ws.Cells.Item(134, 14).Value = -15584.0001  # N134: -14981.454 -> -15584.0001

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(22, 8).Value = 941.7143  # H22: 1100.1052 -> 941.7143
$ws.Cells.Item(22, 9).Value = 366.27274  # I22: 384.9 -> 366.27274
$ws.Cells.Item(22, 10).Value = 1574.7  # J22: 1894.7778 -> 1574.7
$ws.Cells.Item(22, 11).Value = 366.27274  # K22: 384.9 -> 366.27274
$ws.Cells.Item(22, 12).Value = 1574.7  # L22: 1894.7778 -> 1574.7
$ws.Cells.Item(22, 13).Value = -16.27274  # M22: -34.89999999999998 -> -16.27274
$ws.Cells.Item(22, 14).Value = -2274.7  # N22: -2594.7778 -> -2274.7
$ws.Cells.Item(31, 8).Value = 1378.4517  # H31: 1343.7576 -> 1378.4517
$ws.Cells.Item(31, 9).Value = 1021.86957  # I31: 991.5 -> 1021.86957
$ws.Cells.Item(31, 10).Value = 2403.625  # J31: 2283.111 -> 2403.625
$ws.Cells.Item(31, 11).Value = 1021.86957  # K31: 991.5 -> 1021.86957
$ws.Cells.Item(31, 12).Value = 2403.625  # L31: 2283.111 -> 2403.625
$ws.Cells.Item(31, 13).Value = -726.86957  # M31: -696.5 -> -726.86957
$ws.Cells.Item(31, 14).Value = -2993.625  # N31: -2873.111 -> -2993.625
$ws.Cells.Item(34, 8).Value = 1378.4517  # H34: 1343.7576 -> 1378.4517
$ws.Cells.Item(34, 9).Value = 1021.86957  # I34: 991.5 -> 1021.86957
$ws.Cells.Item(34, 10).Value = 2403.625  # J34: 2283.111 -> 2403.625
$ws.Cells.Item(34, 11).Value = 1021.86957  # K34: 991.5 -> 1021.86957
$ws.Cells.Item(34, 12).Value = 2403.625  # L34: 2283.111 -> 2403.625
$ws.Cells.Item(34, 13).Value = -819.86957  # M34: -789.5 -> -819.86957
$ws.Cells.Item(34, 14).Value = -2807.625  # N34: -2687.111 -> -2807.625
$ws.Cells.Item(99, 8).Value = 6205.7856  # H99: 5913.1333 -> 6205.7856
$ws.Cells.Item(99, 10).Value = 7858.6  # J99: 6851.5 -> 7858.6
$ws.Cells.Item(99, 12).Value = 7858.6  # L99: 6851.5 -> 7858.6
$ws.Cells.Item(99, 14).Value = -10854.6  # N99: -9847.5 -> -10854.6
$ws.Cells.Item(105, 8).Value = 14137.857  # H105: 12476.25 -> 14137.857
$ws.Cells.Item(105, 10).Value = 1699.5  # J105: 1272.25 -> 1699.5
$ws.Cells.Item(105, 12).Value = 1699.5  # L105: 1272.25 -> 1699.5
$ws.Cells.Item(105, 14).Value = -5193.5  # N105: -4766.25 -> -5193.5
$ws.Cells.Item(107, 8).Value = 557.6667  # H107: 520.88464 -> 557.6667
$ws.Cells.Item(107, 9).Value = 655.25  # I107: 463.33334 -> 655.25
$ws.Cells.Item(107, 11).Value = 655.25  # K107: 463.33334 -> 655.25
$ws.Cells.Item(107, 13).Value = 1264.75  # M107: 1456.66666 -> 1264.75
$ws.Cells.Item(126, 8).Value = 6205.7856  # H126: 5913.1333 -> 6205.7856
$ws.Cells.Item(126, 10).Value = 7858.6  # J126: 6851.5 -> 7858.6
$ws.Cells.Item(126, 12).Value = 23575.8  # L126: 20554.5 -> 23575.8
$ws.Cells.Item(126, 14).Value = -28515.8  # N126: -25494.5 -> -28515.8
$ws.Cells.Item(134, 8).Value = 1538  # H134: 1610.6923 -> 1538
$ws.Cells.Item(134, 9).Value = 1066.25  # I134: 1202.6 -> 1066.25
$ws.Cells.Item(134, 10).Value = 2953.25  # J134: 2971 -> 2953.25
$ws.Cells.Item(134, 11).Value = 3198.75  # K134: 3607.8 -> 3198.75
$ws.Cells.Item(134, 12).Value = 8859.75  # L134: 8913 -> 8859.75
$ws.Cells.Item(134, 13).Value = -663.75  # M134: -1072.8 -> -663.75
$ws.Cells.Item(134, 14).Value = -13929.75  # N134: -13983 -> -13929.75

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(11, 8).Value = 223  # H11: 175.25 -> 223
$ws.Cells.Item(11, 9).Value = 243.85715  # I11: 133.6 -> 243.85715
$ws.Cells.Item(11, 10).Value = 150  # J11: 244.66667 -> 150
$ws.Cells.Item(11, 11).Value = 731.5714499999999  # K11: 400.8 -> 731.5714499999999
$ws.Cells.Item(11, 12).Value = 450  # L11: 734.00001 -> 450
$ws.Cells.Item(11, 13).Value = -591.5714499999999  # M11: -260.8 -> -591.5714499999999
$ws.Cells.Item(11, 14).Value = -730  # N11: -1014.00001 -> -730
$ws.Cells.Item(75, 8).Value = 0  # H75: 1341.7142 -> 0
$ws.Cells.Item(75, 9).Value = 0  # I75: 1274.25 -> 0
$ws.Cells.Item(75, 10).Value = 0  # J75: 1431.6666 -> 0
$ws.Cells.Item(75, 11).Value = 0  # K75: 3822.75 -> 0
$ws.Cells.Item(75, 12).Value = 0  # L75: 4294.9998 -> 0
$ws.Cells.Item(75, 13).ClearContents()  # M75 was -2824.75
$ws.Cells.Item(75, 14).ClearContents()  # N75 was -6290.9998
$ws.Cells.Item(78, 8).Value = 0  # H78: 1341.7142 -> 0
$ws.Cells.Item(78, 9).Value = 0  # I78: 1274.25 -> 0
$ws.Cells.Item(78, 10).Value = 0  # J78: 1431.6666 -> 0
$ws.Cells.Item(78, 11).Value = 0  # K78: 11468.25 -> 0
$ws.Cells.Item(78, 12).Value = 0  # L78: 12884.9994 -> 0
$ws.Cells.Item(78, 13).ClearContents()  # M78 was -6476.25
$ws.Cells.Item(78, 14).ClearContents()  # N78 was -22868.9994
$ws.Cells.Item(113, 8).Value = 2102.4  # H113: 2181.7144 -> 2102.4
$ws.Cells.Item(113, 10).Value = 1628.4166  # J113: 1686.2727 -> 1628.4166
$ws.Cells.Item(113, 12).Value = 4885.2498  # L113: 5058.8181 -> 4885.2498
$ws.Cells.Item(113, 14).Value = -9225.2498  # N113: -9398.8181 -> -9225.2498
$ws.Cells.Item(129, 8).Value = 868.4286  # H129: 827.5 -> 868.4286
$ws.Cells.Item(129, 9).Value = 938.6667  # I129: 827.5 -> 938.6667
$ws.Cells.Item(129, 10).Value = 447  # J129: 0 -> 447
$ws.Cells.Item(129, 11).Value = 2816.0001  # K129: 2482.5 -> 2816.0001
$ws.Cells.Item(129, 12).Value = 1341  # L129: 0 -> 1341
$ws.Cells.Item(129, 13).Value = 2183.9999  # M129: 2517.5 -> 2183.9999
$ws.Cells.Item(129, 14).Value = -11341  # N129: None -> -11341
$ws.Cells.Item(131, 8).Value = 1846.8868  # H131: 2009.8864 -> 1846.8868
$ws.Cells.Item(131, 9).Value = 1008.3333  # I131: 1015 -> 1008.3333
$ws.Cells.Item(131, 10).Value = 1897.2  # J131: 2057.262 -> 1897.2
$ws.Cells.Item(131, 11).Value = 3024.9999  # K131: 3045 -> 3024.9999
$ws.Cells.Item(131, 12).Value = 5691.6  # L131: 6171.786 -> 5691.6
$ws.Cells.Item(131, 13).Value = 2015.0001  # M131: 1995 -> 2015.0001
$ws.Cells.Item(131, 14).Value = -15771.6  # N131: -16251.786 -> -15771.6
$ws.Cells.Item(140, 8).Value = 2746.3333  # H140: 3313.5 -> 2746.3333
$ws.Cells.Item(140, 9).Value = 2145  # I140: 2522.5 -> 2145
$ws.Cells.Item(140, 10).Value = 4400  # J140: 4500 -> 4400
$ws.Cells.Item(140, 11).Value = 6435  # K140: 7567.5 -> 6435
$ws.Cells.Item(140, 12).Value = 13200  # L140: 13500 -> 13200
$ws.Cells.Item(140, 13).Value = -1255  # M140: -2387.5 -> -1255
$ws.Cells.Item(140, 14).Value = -23560  # N140: -23860 -> -23560

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(18, 8).Value = 0  # H18: 2999 -> 0
$ws.Cells.Item(18, 9).Value = 0  # I18: 2999 -> 0
$ws.Cells.Item(18, 11).Value = 0  # K18: 2999 -> 0
$ws.Cells.Item(18, 13).ClearContents()  # M18 was -2706
$ws.Cells.Item(80, 8).Value = 7808.5  # H80: 12614 -> 7808.5
$ws.Cells.Item(80, 9).Value = 12611  # I80: 22222 -> 12611
$ws.Cells.Item(80, 11).Value = 12611  # K80: 22222 -> 12611
$ws.Cells.Item(80, 13).Value = -11613  # M80: -21224 -> -11613
$ws.Cells.Item(83, 8).Value = 7808.5  # H83: 12614 -> 7808.5
$ws.Cells.Item(83, 9).Value = 12611  # I83: 22222 -> 12611
$ws.Cells.Item(83, 11).Value = 63055  # K83: 111110 -> 63055
$ws.Cells.Item(83, 13).Value = -58063  # M83: -106118 -> -58063
$ws.Cells.Item(122, 8).Value = 2260.6562  # H122: 2261.5 -> 2260.6562
$ws.Cells.Item(122, 9).Value = 1787.5769  # I122: 1788.6154 -> 1787.5769
$ws.Cells.Item(122, 11).Value = 5362.7307  # K122: 5365.8462 -> 5362.7307
$ws.Cells.Item(122, 13).Value = -2912.7307  # M122: -2915.8462 -> -2912.7307
$ws.Cells.Item(132, 8).Value = 2297.55  # H132: 3347.5 -> 2297.55
$ws.Cells.Item(132, 9).Value = 2041.6  # I132: 2486.25 -> 2041.6
$ws.Cells.Item(132, 10).Value = 2553.5  # J132: 5070 -> 2553.5
$ws.Cells.Item(132, 11).Value = 6124.799999999999  # K132: 7458.75 -> 6124.799999999999
$ws.Cells.Item(132, 12).Value = 7660.5  # L132: 15210 -> 7660.5
$ws.Cells.Item(132, 13).Value = -3594.799999999999  # M132: -4928.75 -> -3594.799999999999
$ws.Cells.Item(132, 14).Value = -12720.5  # N132: -20270 -> -12720.5

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 2823.4546  # H7: 2793.4707 -> 2823.4546
$ws.Cells.Item(7, 9).Value = 2282.4  # I7: 2252.5 -> 2282.4
$ws.Cells.Item(7, 11).Value = 2282.4  # K7: 2252.5 -> 2282.4
$ws.Cells.Item(7, 13).Value = -2170.4  # M7: -2140.5 -> -2170.4
$ws.Cells.Item(14, 8).Value = 2602  # H14: 3004 -> 2602
$ws.Cells.Item(14, 9).Value = 2602  # I14: 3004 -> 2602
$ws.Cells.Item(14, 11).Value = 2602  # K14: 3004 -> 2602
$ws.Cells.Item(14, 13).Value = -2430  # M14: -2832 -> -2430
$ws.Cells.Item(42, 8).Value = 22400  # H42: 21800.666 -> 22400
$ws.Cells.Item(42, 9).Value = 0  # I42: 21451 -> 0
$ws.Cells.Item(42, 10).Value = 22400  # J42: 22500 -> 22400
$ws.Cells.Item(42, 11).Value = 0  # K42: 21451 -> 0
$ws.Cells.Item(42, 12).Value = 22400  # L42: 22500 -> 22400
$ws.Cells.Item(42, 13).ClearContents()  # M42 was -20888
$ws.Cells.Item(42, 14).Value = -23526  # N42: -23626 -> -23526
$ws.Cells.Item(43, 8).Value = 20056  # H43: 20295.334 -> 20056
$ws.Cells.Item(43, 10).Value = 20056  # J43: 20295.334 -> 20056
$ws.Cells.Item(43, 12).Value = 20056  # L43: 20295.334 -> 20056
$ws.Cells.Item(43, 14).Value = -20442  # N43: -20681.334 -> -20442
$ws.Cells.Item(49, 8).Value = 22400  # H49: 21800.666 -> 22400
$ws.Cells.Item(49, 9).Value = 0  # I49: 21451 -> 0
$ws.Cells.Item(49, 10).Value = 22400  # J49: 22500 -> 22400
$ws.Cells.Item(49, 11).Value = 0  # K49: 21451 -> 0
$ws.Cells.Item(49, 12).Value = 22400  # L49: 22500 -> 22400
$ws.Cells.Item(49, 13).ClearContents()  # M49 was -21304
$ws.Cells.Item(49, 14).Value = -22694  # N49: -22794 -> -22694
$ws.Cells.Item(61, 8).Value = 3072.2856  # H61: 3358.1428 -> 3072.2856
$ws.Cells.Item(61, 10).Value = 3252.5  # J61: 3752.75 -> 3252.5
$ws.Cells.Item(61, 12).Value = 3252.5  # L61: 3752.75 -> 3252.5
$ws.Cells.Item(61, 14).Value = -3656.5  # N61: -4156.75 -> -3656.5
$ws.Cells.Item(113, 8).Value = 3072.2856  # H113: 3358.1428 -> 3072.2856
$ws.Cells.Item(113, 10).Value = 3252.5  # J113: 3752.75 -> 3252.5
$ws.Cells.Item(113, 12).Value = 3252.5  # L113: 3752.75 -> 3252.5
$ws.Cells.Item(113, 14).Value = -7592.5  # N113: -8092.75 -> -7592.5
$ws.Cells.Item(122, 8).Value = 7897.6  # H122: 8102.8423 -> 7897.6
$ws.Cells.Item(122, 9).Value = 5746.2144  # I122: 5880.6924 -> 5746.2144
$ws.Cells.Item(122, 11).Value = 17238.6432  # K122: 17642.0772 -> 17238.6432
$ws.Cells.Item(122, 13).Value = -14788.6432  # M122: -15192.0772 -> -14788.6432
$ws.Cells.Item(126, 8).Value = 2823.4546  # H126: 2793.4707 -> 2823.4546
$ws.Cells.Item(126, 9).Value = 2282.4  # I126: 2252.5 -> 2282.4
$ws.Cells.Item(126, 11).Value = 6847.200000000001  # K126: 6757.5 -> 6847.200000000001
$ws.Cells.Item(126, 13).Value = -4377.200000000001  # M126: -4287.5 -> -4377.200000000001
$ws.Cells.Item(136, 8).Value = 3048.9473  # H136: 3740 -> 3048.9473
$ws.Cells.Item(136, 9).Value = 3135  # I136: 3880 -> 3135
$ws.Cells.Item(136, 11).Value = 9405  # K136: 11640 -> 9405
$ws.Cells.Item(136, 13).Value = -6855  # M136: -9090 -> -6855

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(49, 8).Value = 20000  # H49: 0 -> 20000
$ws.Cells.Item(49, 9).Value = 20000  # I49: 0 -> 20000
$ws.Cells.Item(49, 11).Value = 20000  # K49: 0 -> 20000
$ws.Cells.Item(49, 13).Value = -19770  # M49: None -> -19770
$ws.Cells.Item(74, 8).Value = 29134.8  # H74: 29134.6 -> 29134.8
$ws.Cells.Item(74, 9).Value = 0  # I74: 23889 -> 0
$ws.Cells.Item(74, 10).Value = 29134.8  # J74: 30446 -> 29134.8
$ws.Cells.Item(74, 11).Value = 0  # K74: 23889 -> 0
$ws.Cells.Item(74, 12).Value = 29134.8  # L74: 30446 -> 29134.8
$ws.Cells.Item(74, 13).ClearContents()  # M74 was -22953
$ws.Cells.Item(74, 14).Value = -31006.8  # N74: -32318 -> -31006.8
$ws.Cells.Item(77, 8).Value = 29134.8  # H77: 29134.6 -> 29134.8
$ws.Cells.Item(77, 9).Value = 0  # I77: 23889 -> 0
$ws.Cells.Item(77, 10).Value = 29134.8  # J77: 30446 -> 29134.8
$ws.Cells.Item(77, 11).Value = 0  # K77: 71667 -> 0
$ws.Cells.Item(77, 12).Value = 87404.39999999999  # L77: 91338 -> 87404.39999999999
$ws.Cells.Item(77, 13).ClearContents()  # M77 was -66987
$ws.Cells.Item(77, 14).Value = -96764.39999999999  # N77: -100698 -> -96764.39999999999
$ws.Cells.Item(113, 8).Value = 1574.25  # H113: 1727.7142 -> 1574.25
$ws.Cells.Item(113, 10).Value = 1674.625  # J113: 2066.1667 -> 1674.625
$ws.Cells.Item(113, 12).Value = 5023.875  # L113: 6198.500100000001 -> 5023.875
$ws.Cells.Item(113, 14).Value = -9363.875  # N113: -10538.5001 -> -9363.875
$ws.Cells.Item(122, 8).Value = 4130.8335  # H122: 4075.4324 -> 4130.8335
$ws.Cells.Item(122, 9).Value = 4679.4585  # I122: 4302.2593 -> 4679.4585
$ws.Cells.Item(122, 10).Value = 3033.5833  # J122: 3463 -> 3033.5833
$ws.Cells.Item(122, 11).Value = 14038.3755  # K122: 12906.7779 -> 14038.3755
$ws.Cells.Item(122, 12).Value = 9100.749899999999  # L122: 10389 -> 9100.749899999999
$ws.Cells.Item(122, 13).Value = -11588.3755  # M122: -10456.7779 -> -11588.3755
$ws.Cells.Item(122, 14).Value = -14000.7499  # N122: -15289 -> -14000.7499
$ws.Cells.Item(126, 8).Value = 2327.9333  # H126: 2328.0667 -> 2327.9333
$ws.Cells.Item(126, 10).Value = 2448.125  # J126: 2448.375 -> 2448.125
$ws.Cells.Item(126, 12).Value = 7344.375  # L126: 7345.125 -> 7344.375
$ws.Cells.Item(126, 14).Value = -12284.375  # N126: -12285.125 -> -12284.375
$ws.Cells.Item(132, 8).Value = 1565.4546  # H132: 1668.9445 -> 1565.4546
$ws.Cells.Item(132, 9).Value = 1482.579  # I132: 1584.6666 -> 1482.579
$ws.Cells.Item(132, 11).Value = 4447.737  # K132: 4753.9998 -> 4447.737
$ws.Cells.Item(132, 13).Value = -1917.737  # M132: -2223.9998 -> -1917.737
$ws.Cells.Item(136, 8).Value = 6441.75  # H136: 7555.4 -> 6441.75
$ws.Cells.Item(136, 9).Value = 4667.15  # I136: 5615.5625 -> 4667.15
$ws.Cells.Item(136, 11).Value = 14001.45  # K136: 16846.6875 -> 14001.45
$ws.Cells.Item(136, 13).Value = -11451.45  # M136: -14296.6875 -> -11451.45
